$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$shp = $ws.Shapes.AddShape(9, 100, 100, 200, 80)
$shp.TextFrame.Characters.Text = "Hello"
Write-Output $shp.Name
Write-Output $shp.Type
